# Auto-generated PowerShell COM-interop script to apply the diff to Bahamut_Profits sheets
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across the ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 15750
$ws.Range("I59").Value = 2000
$ws.Range("J59").Value = 20333.334
$ws.Range("K59").Value = 6000
$ws.Range("L59").Value = 61000.00199999999
$ws.Range("M59").Value = -5443
$ws.Range("N59").Value = -62114.00199999999
$ws.Range("H62").Value = 134418.22
$ws.Range("I62").Value = 134418.22
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 134418.22
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -133794.22
$ws.Range("H65").Value = 134418.22
$ws.Range("I65").Value = 134418.22
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 672091.1
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -668971.1
$ws.Range("H106").Value = 2481.3572
$ws.Range("I106").Value = 2353.25
$ws.Range("K106").Value = 2353.25
$ws.Range("M106").Value = -1722.25
$ws.Range("H125").Value = 17997.75
$ws.Range("I125").Value = 1001
$ws.Range("J125").Value = 23663.334
$ws.Range("K125").Value = 9009
$ws.Range("L125").Value = 212970.006
$ws.Range("M125").Value = -6549
$ws.Range("N125").Value = -217890.006
$ws.Range("H129").Value = 1950141.1
$ws.Range("J129").Value = 3087452.2
$ws.Range("L129").Value = 9262356.600000001
$ws.Range("N129").Value = -9272356.600000001
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1568
$ws.Range("I61").Value = 822.4
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 822.4
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -610.4
$ws.Range("N61").Value = -2924
$ws.Range("H122").Value = 1765
$ws.Range("I122").Value = 1588.5714
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4765.7142
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2315.7142
$ws.Range("N122").Value = -13900
$ws.Range("H136").Value = 1568
$ws.Range("I136").Value = 822.4
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 2467.2
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = 82.80000000000018
$ws.Range("N136").Value = -12600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H54").Value = 1474.5
$ws.Range("I54").Value = 1474.5
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1474.5
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -990.5
$ws.Range("H117").Value = 39742
$ws.Range("J117").Value = 39742
$ws.Range("L117").Value = 39742
$ws.Range("N117").Value = -48920
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("N52").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("N121").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 55000
$ws.Range("J119").Value = 55000
$ws.Range("L119").Value = 55000
$ws.Range("N119").Value = -64676
$ws.Range("H123").Value = 51893
$ws.Range("J123").Value = 51893
$ws.Range("L123").Value = 51893
$ws.Range("N123").Value = -61693

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 415
$ws.Range("I15").Value = 80
$ws.Range("J15").Value = 750
$ws.Range("K15").Value = 240
$ws.Range("L15").Value = 2250
$ws.Range("M15").Value = -100
$ws.Range("N15").Value = -2530
$ws.Range("H38").Value = 180.125
$ws.Range("I38").Value = 60.25
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 180.75
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = 166.25
$ws.Range("N38").Value = -1594
$ws.Range("H42").Value = 2322.3333
$ws.Range("I42").Value = 500
$ws.Range("J42").Value = 3233.5
$ws.Range("K42").Value = 1500
$ws.Range("L42").Value = 9700.5
$ws.Range("M42").Value = -966
$ws.Range("N42").Value = -10768.5
$ws.Range("H107").Value = 1004.7273
$ws.Range("I107").Value = 358.66666
$ws.Range("J107").Value = 1148.2963
$ws.Range("K107").Value = 1075.99998
$ws.Range("L107").Value = 3444.8889
$ws.Range("M107").Value = 844.0000199999999
$ws.Range("N107").Value = -7284.8889
$ws.Range("H131").Value = 22427.55
$ws.Range("J131").Value = 2034.5238
$ws.Range("L131").Value = 6103.5714
$ws.Range("N131").Value = -16183.5714

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5226
$ws.Range("H7").Value = 1366.3334
$ws.Range("I7").Value = 1449.5
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 1449.5
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = -1337.5
$ws.Range("N7").Value = -1424
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("K10").Value = 1000
$ws.Range("M10").Value = -860
$ws.Range("H31").Value = 1060.7778
$ws.Range("J31").Value = 1043.4
$ws.Range("L31").Value = 1043.4
$ws.Range("N31").Value = -1539.4
$ws.Range("H40").Value = 1139723.8
$ws.Range("I40").Value = 1139723.8
$ws.Range("K40").Value = 1139723.8
$ws.Range("M40").Value = -1139587.8
$ws.Range("H43").Value = 1000000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H53").Value = 36142.855
$ws.Range("J53").Value = 36333.332
$ws.Range("L53").Value = 36333.332
$ws.Range("N53").Value = -37369.332
$ws.Range("H56").Value = 6000
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H126").Value = 1366.3334
$ws.Range("I126").Value = 1449.5
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 4348.5
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -1878.5
$ws.Range("N126").Value = -8540
$ws.Range("N43").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("N121").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 70017
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 70017
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 70017
$ws.Range("N21").Value = -70487
$ws.Range("H23").Value = 572.25
$ws.Range("I23").Value = 100
$ws.Range("K23").Value = 100
$ws.Range("M23").Value = 129
$ws.Range("H28").Value = 43182.6
$ws.Range("J28").Value = 43182.6
$ws.Range("L28").Value = 43182.6
$ws.Range("N28").Value = -43878.6
$ws.Range("H30").Value = 35254.75
$ws.Range("I30").Value = 1009
$ws.Range("J30").Value = 46670
$ws.Range("K30").Value = 1009
$ws.Range("L30").Value = 46670
$ws.Range("M30").Value = -902
$ws.Range("N30").Value = -46884
$ws.Range("H33").Value = 5000
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5500
$ws.Range("H35").Value = 70017
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70017
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70017
$ws.Range("N35").Value = -70597
$ws.Range("H36").Value = 5000
$ws.Range("J36").Value = 5000
$ws.Range("L36").Value = 5000
$ws.Range("N36").Value = -5500
$ws.Range("H40").Value = 3000
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3298
$ws.Range("H51").Value = 12000
$ws.Range("I51").Value = 500
$ws.Range("K51").Value = 500
$ws.Range("M51").Value = 10
$ws.Range("H53").Value = 4833.3335
$ws.Range("J53").Value = 4833.3335
$ws.Range("L53").Value = 4833.3335
$ws.Range("N53").Value = -6047.3335
$ws.Range("H55").Value = 8333
$ws.Range("J55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("N55").Value = -8554
$ws.Range("M21").ClearContents()
$ws.Range("M35").ClearContents()
